$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.060748854419377
$ws.Cells.Item(2, 4).Value = 1.06224691156679
$ws.Cells.Item(2, 5).Value = 1.073177383261232
$ws.Cells.Item(2, 6).Value = 1.077839669611633
$ws.Cells.Item(2, 9).Value = 1.04509808180284
$ws.Cells.Item(2, 10).Value = 1.065728096918557
$ws.Cells.Item(2, 11).Value = 1.064968777257193
$ws.Cells.Item(2, 12).Value = 1.075869935488028
$ws.Cells.Item(2, 13).Value = 1.080519913913896
$ws.Cells.Item(2, 14).Value = 1.025499034381038

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.06214398129997
$ws.Cells.Item(3, 4).Value = 1.063326797455164
$ws.Cells.Item(3, 5).Value = 1.07451834757054
$ws.Cells.Item(3, 6).Value = 1.079159149565894
$ws.Cells.Item(3, 9).Value = 1.045434359459927
$ws.Cells.Item(3, 10).Value = 1.066774868218241
$ws.Cells.Item(3, 11).Value = 1.065862775052812
$ws.Cells.Item(3, 12).Value = 1.077026495594344
$ws.Cells.Item(3, 13).Value = 1.08165593958843
$ws.Cells.Item(3, 14).Value = 1.025858558975796

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.06304591947832
$ws.Cells.Item(4, 4).Value = 1.064024558777395
$ws.Cells.Item(4, 5).Value = 1.075385552562822
$ws.Cells.Item(4, 6).Value = 1.080012380850949
$ws.Cells.Item(4, 9).Value = 1.045649826885242
$ws.Cells.Item(4, 10).Value = 1.067450927782097
$ws.Cells.Item(4, 11).Value = 1.066439658675902
$ws.Cells.Item(4, 12).Value = 1.077773832541324
$ws.Cells.Item(4, 13).Value = 1.082389911639734
$ws.Cells.Item(4, 14).Value = 1.026090497576959

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.063424906922157
$ws.Cells.Item(5, 4).Value = 1.064317661703892
$ws.Cells.Item(5, 5).Value = 1.075750012378138
$ws.Cells.Item(5, 6).Value = 1.080370948742065
$ws.Cells.Item(5, 9).Value = 1.0457399013949
$ws.Cells.Item(5, 10).Value = 1.067734842066924
$ws.Cells.Item(5, 11).Value = 1.066681801876509
$ws.Cells.Item(5, 12).Value = 1.078087768702679
$ws.Cells.Item(5, 13).Value = 1.082698210450356
$ws.Cells.Item(5, 14).Value = 1.026187838523082

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.063488529730098
$ws.Cells.Item(6, 4).Value = 1.064366861219718
$ws.Cells.Item(6, 5).Value = 1.075811200259109
$ws.Cells.Item(6, 6).Value = 1.080431146319158
$ws.Cells.Item(6, 9).Value = 1.045754995545265
$ws.Cells.Item(6, 10).Value = 1.067782494958414
$ws.Cells.Item(6, 11).Value = 1.066722436655787
$ws.Cells.Item(6, 12).Value = 1.078140465739772
$ws.Cells.Item(6, 13).Value = 1.082749959848061
$ws.Cells.Item(6, 14).Value = 1.026204172789195

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.063050984257505
$ws.Cells.Item(7, 4).Value = 1.06402847615669
$ws.Cells.Item(7, 5).Value = 1.075390422932652
$ws.Cells.Item(7, 6).Value = 1.080017172561091
$ws.Cells.Item(7, 9).Value = 1.045651032459928
$ws.Cells.Item(7, 10).Value = 1.067454722638002
$ws.Cells.Item(7, 11).Value = 1.066442895689892
$ws.Cells.Item(7, 12).Value = 1.077778028327322
$ws.Cells.Item(7, 13).Value = 1.082394032173015
$ws.Cells.Item(7, 14).Value = 1.026091798902605

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.06122051300592
$ws.Cells.Item(8, 4).Value = 1.062612072175356
$ws.Cells.Item(8, 5).Value = 1.073630671887644
$ws.Cells.Item(8, 6).Value = 1.078285712285573
$ws.Cells.Item(8, 9).Value = 1.045212169450198
$ws.Cells.Item(8, 10).Value = 1.066082123396977
$ws.Cells.Item(8, 11).Value = 1.065271239370365
$ws.Cells.Item(8, 12).Value = 1.076261016722597
$ws.Cells.Item(8, 13).Value = 1.080904071456911
$ws.Cells.Item(8, 14).Value = 1.02562068251487

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.057988634559585
$ws.Cells.Item(9, 4).Value = 1.060108426396481
$ws.Cells.Item(9, 5).Value = 1.070525852604753
$ws.Cells.Item(9, 6).Value = 1.075230201729793
$ws.Cells.Item(9, 9).Value = 1.04442249697978
$ws.Cells.Item(9, 10).Value = 1.06365354468553
$ws.Cells.Item(9, 11).Value = 1.063194314277673
$ws.Cells.Item(9, 12).Value = 1.073579768605794
$ws.Cells.Item(9, 13).Value = 1.078269901074075
$ws.Cells.Item(9, 14).Value = 1.024785124336948

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.055829480541024
$ws.Cells.Item(10, 4).Value = 1.058433933930044
$ws.Cells.Item(10, 5).Value = 1.068453100305887
$ws.Cells.Item(10, 6).Value = 1.073189966876898
$ws.Cells.Item(10, 9).Value = 1.043884988596998
$ws.Cells.Item(10, 10).Value = 1.062027641261956
$ws.Cells.Item(10, 11).Value = 1.061801248081167
$ws.Cells.Item(10, 12).Value = 1.071786618724033
$ws.Cells.Item(10, 13).Value = 1.076507748414759
$ws.Cells.Item(10, 14).Value = 1.024224395479185

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.054893392912321
$ws.Cells.Item(11, 4).Value = 1.057707540215316
$ws.Cells.Item(11, 5).Value = 1.067554839709766
$ws.Cells.Item(11, 6).Value = 1.072305704708891
$ws.Cells.Item(11, 9).Value = 1.043649601183944
$ws.Cells.Item(11, 10).Value = 1.06132193566563
$ws.Cells.Item(11, 11).Value = 1.061195992735054
$ws.Cells.Item(11, 12).Value = 1.071008775243802
$ws.Cells.Item(11, 13).Value = 1.075743236880316
$ws.Cells.Item(11, 14).Value = 1.023980703581787

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.054545508236429
$ws.Cells.Item(12, 4).Value = 1.057437522238343
$ws.Cells.Item(12, 5).Value = 1.067221068977438
$ws.Cells.Item(12, 6).Value = 1.071977121507901
$ws.Cells.Item(12, 9).Value = 1.043561769227192
$ws.Cells.Item(12, 10).Value = 1.061059548880696
$ws.Cells.Item(12, 11).Value = 1.060970862795412
$ws.Cells.Item(12, 12).Value = 1.070719635191966
$ws.Cells.Item(12, 13).Value = 1.0754590355298
$ws.Cells.Item(12, 14).Value = 1.023890050134534

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.054620138921854
$ws.Cells.Item(13, 4).Value = 1.057495451255312
$ws.Cells.Item(13, 5).Value = 1.067292669293174
$ws.Cells.Item(13, 6).Value = 1.072047609625463
$ws.Cells.Item(13, 9).Value = 1.043580627542177
$ws.Cells.Item(13, 10).Value = 1.061115843401236
$ws.Cells.Item(13, 11).Value = 1.061019168066852
$ws.Cells.Item(13, 12).Value = 1.070781666484587
$ws.Cells.Item(13, 13).Value = 1.075520008057102
$ws.Cells.Item(13, 14).Value = 1.023909501756353

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.05486464036327
$ws.Cells.Item(14, 4).Value = 1.057685224607216
$ws.Cells.Item(14, 5).Value = 1.067527252526479
$ws.Cells.Item(14, 6).Value = 1.07227854657448
$ws.Cells.Item(14, 9).Value = 1.043642349107148
$ws.Cells.Item(14, 10).Value = 1.061300251944499
$ws.Cells.Item(14, 11).Value = 1.061177389791199
$ws.Cells.Item(14, 12).Value = 1.070984879234877
$ws.Cells.Item(14, 13).Value = 1.075719749377984
$ws.Cells.Item(14, 14).Value = 1.023973212909893

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.055015261737649
$ws.Cells.Item(15, 4).Value = 1.05780212322452
$ws.Cells.Item(15, 5).Value = 1.067671771383926
$ws.Cells.Item(15, 6).Value = 1.072420817236998
$ws.Cells.Item(15, 9).Value = 1.043680324936633
$ws.Cells.Item(15, 10).Value = 1.061413838048212
$ws.Cells.Item(15, 11).Value = 1.061274834084554
$ws.Cells.Item(15, 12).Value = 1.071110056783038
$ws.Cells.Item(15, 13).Value = 1.075842786295915
$ws.Cells.Item(15, 14).Value = 1.024012449469321

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.055891580664584
$ws.Cells.Item(16, 4).Value = 1.058482114016462
$ws.Cells.Item(16, 5).Value = 1.068512698794529
$ws.Cells.Item(16, 6).Value = 1.073248634627683
$ws.Cells.Item(16, 9).Value = 1.043900554673004
$ws.Cells.Item(16, 10).Value = 1.062074440900766
$ws.Cells.Item(16, 11).Value = 1.061841373436791
$ws.Cells.Item(16, 12).Value = 1.071838211814117
$ws.Cells.Item(16, 13).Value = 1.076558454823188
$ws.Cells.Item(16, 14).Value = 1.024240549589027

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.056440957793888
$ws.Cells.Item(17, 4).Value = 1.058908296367651
$ws.Cells.Item(17, 5).Value = 1.069039987460083
$ws.Cells.Item(17, 6).Value = 1.073767677945398
$ws.Cells.Item(17, 9).Value = 1.044037990342848
$ws.Cells.Item(17, 10).Value = 1.062488367298101
$ws.Cells.Item(17, 11).Value = 1.062196197778074
$ws.Cells.Item(17, 12).Value = 1.072294586945261
$ws.Cells.Item(17, 13).Value = 1.077006973497808
$ws.Cells.Item(17, 14).Value = 1.024383390861576

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.056761288681171
$ws.Cells.Item(18, 4).Value = 1.059156753489331
$ws.Cells.Item(18, 5).Value = 1.069347474446021
$ws.Cells.Item(18, 6).Value = 1.074070347650379
$ws.Cells.Item(18, 9).Value = 1.044117899364046
$ws.Cells.Item(18, 10).Value = 1.06272964187315
$ws.Cells.Item(18, 11).Value = 1.062402963394917
$ws.Cells.Item(18, 12).Value = 1.072560648272031
$ws.Cells.Item(18, 13).Value = 1.077268443586997
$ws.Cells.Item(18, 14).Value = 1.024466621726196

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.05687049449894
$ws.Cells.Item(19, 4).Value = 1.059241449371782
$ws.Cells.Item(19, 5).Value = 1.06945230756653
$ws.Cells.Item(19, 6).Value = 1.074173536855376
$ws.Cells.Item(19, 9).Value = 1.044145103089666
$ws.Cells.Item(19, 10).Value = 1.062811882966138
$ws.Cells.Item(19, 11).Value = 1.062473431743301
$ws.Cells.Item(19, 12).Value = 1.072651345615001
$ws.Cells.Item(19, 13).Value = 1.077357573956394
$ws.Cells.Item(19, 14).Value = 1.024494986747215

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.056382026403594
$ws.Cells.Item(20, 4).Value = 1.058862584284209
$ws.Cells.Item(20, 5).Value = 1.068983421815642
$ws.Cells.Item(20, 6).Value = 1.073711997765438
$ws.Cells.Item(20, 9).Value = 1.044023271171343
$ws.Cells.Item(20, 10).Value = 1.062443973668245
$ws.Cells.Item(20, 11).Value = 1.062158148929823
$ws.Cells.Item(20, 12).Value = 1.072245636164399
$ws.Cells.Item(20, 13).Value = 1.076958866539072
$ws.Cells.Item(20, 14).Value = 1.024368074253674

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.054792645800021
$ws.Cells.Item(21, 4).Value = 1.057629346717958
$ws.Cells.Item(21, 5).Value = 1.067458176872713
$ws.Cells.Item(21, 6).Value = 1.072210544993294
$ws.Cells.Item(21, 9).Value = 1.043624184660208
$ws.Cells.Item(21, 10).Value = 1.061245955333112
$ws.Cells.Item(21, 11).Value = 1.061130806053725
$ws.Cells.Item(21, 12).Value = 1.070925044097071
$ws.Cells.Item(21, 13).Value = 1.075660936857657
$ws.Cells.Item(21, 14).Value = 1.023954455312226

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.053792294504388
$ws.Cells.Item(22, 4).Value = 1.056852784975974
$ws.Cells.Item(22, 5).Value = 1.066498516951738
$ws.Cells.Item(22, 6).Value = 1.071265774193293
$ws.Cells.Item(22, 9).Value = 1.043370956229322
$ws.Cells.Item(22, 10).Value = 1.060491228705059
$ws.Cells.Item(22, 11).Value = 1.060483072866789
$ws.Cells.Item(22, 12).Value = 1.070093493111828
$ws.Cells.Item(22, 13).Value = 1.074843557294713
$ws.Cells.Item(22, 14).Value = 1.02369361247188

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.054322700623438
$ws.Cells.Item(23, 4).Value = 1.057264567752611
$ws.Cells.Item(23, 5).Value = 1.067007316698387
$ws.Cells.Item(23, 6).Value = 1.071766687504105
$ws.Cells.Item(23, 9).Value = 1.043505416600109
$ws.Cells.Item(23, 10).Value = 1.060891465702408
$ws.Cells.Item(23, 11).Value = 1.060826620426958
$ws.Cells.Item(23, 12).Value = 1.070534433242217
$ws.Cells.Item(23, 13).Value = 1.07527699216392
$ws.Cells.Item(23, 14).Value = 1.023831964973168

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.056408655306347
$ws.Cells.Item(24, 4).Value = 1.058883240001888
$ws.Cells.Item(24, 5).Value = 1.069008981616419
$ws.Cells.Item(24, 6).Value = 1.073737157487387
$ws.Cells.Item(24, 9).Value = 1.044029922918831
$ws.Cells.Item(24, 10).Value = 1.062464033738042
$ws.Cells.Item(24, 11).Value = 1.062175342176369
$ws.Cells.Item(24, 12).Value = 1.072267755330761
$ws.Cells.Item(24, 13).Value = 1.076980604445431
$ws.Cells.Item(24, 14).Value = 1.024374995435258

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.058824936007658
$ws.Cells.Item(25, 4).Value = 1.060756616718505
$ws.Cells.Item(25, 5).Value = 1.071329011338655
$ws.Cells.Item(25, 6).Value = 1.076020676300149
$ws.Cells.Item(25, 9).Value = 1.044628590173161
$ws.Cells.Item(25, 10).Value = 1.064282582849619
$ws.Cells.Item(25, 11).Value = 1.063732726061513
$ws.Cells.Item(25, 12).Value = 1.074273915445908
$ws.Cells.Item(25, 13).Value = 1.078951946177599
$ws.Cells.Item(25, 14).Value = 1.025001781724772

